# "Bauluz updates and others"
# - fix a data-entry value (B58: 3 -> 4)
# - append 6 new country/quality rows to the Tableau1 table (ER, GQ, LY, SO, ZZ, QF)
# - move the selection/scroll position down to the new last row (A210)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing quality value for row 58 (B58: 3 -> 4)
$ws.Cells.Item(58, 2).Value2 = 4

# The data lives in an Excel Table ("Tableau1"); grow it so the new rows are
# picked up by the table range / autofilter / dimension automatically.
$lo = $ws.ListObjects.Item(1)

$newRows = @(
    @("ER", 0),
    @("GQ", 0),
    @("LY", 0),
    @("SO", 0),
    @("ZZ", 0),
    @("QF", 3)
)

foreach ($entry in $newRows) {
    $newRow = $lo.ListRows.Add()
    $rowIndex = $newRow.Range.Row
    $ws.Cells.Item($rowIndex, 1).Value2 = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value2 = $entry[1]
}

# Bring the view to the newly added tail of the table, like the saved
# workbook (topLeftCell A206 / selection A210).
$excel.ActiveWindow.ScrollRow = 206
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A210").Select()
